$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "deep sleep" power value (duty-cycle data entry). The dependent
# formulas (F10 = D10*E10, F11 = SUM(F4:F10), I11 = F10/F11, D16, D17, D18)
# recalculate automatically from this single input change.
$ws.Range("E10").Value = 59.66

# Reflect the updated view/scroll state: the window was scrolled one column
# right (B3 -> C3) and the active selection moved down one row (E10 -> E11).
$excel.ActiveWindow.ScrollColumn = 3
$excel.ActiveWindow.ScrollRow = 3
$ws.Range("E11").Select()
